$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Date row (A8 label "Date", B8 value)
$ws.Range("B8").Value = "2026-01-28T10:29:57+00:00"

# Description row (A13 label "Description", B13 value) - remove the extra second line
$ws.Range("B13").Value = "Libellé court de l'organisation"

# Context row (A22 label "Context", B22 value)
$ws.Range("B22").Value = "element:https://hl7.fr/ig/fhir/core/StructureDefinition/fr-core-organization#Organization"
